{"js": "// ---------------------------------------------------------------------------\n// Target change (see XML diff in the task description):\n//   word/numbering.xml - for the ten legacy <w:abstractNum> definitions\n//   (abstractNumId=\"1\"..\"10\", the unused \"FFFFFF7C\"..\"FFFFFF89\" nsid block\n//   inherited from the Pandoc/R-Markdown \"reference docx\" template) the\n//   cached <w:tmpl w:val=\"\u2026\"/> correlation GUID is swapped for a new,\n//   unrelated GUID. Nothing else in those definitions (nsid, lvl, numFmt,\n//   indentation, \u2026) changes, and no paragraph in the document body actually\n//   uses any of these numbering definitions (no <w:numPr> anywhere).\n//\n// <w:tmpl> is Word's internal \"list template cache key\". It is not content,\n// not formatting, and not surfaced by any read/write property anywhere in\n// the Word JavaScript API - there is no Range/Paragraph/List/ListLevel (or\n// any other) object whose property maps to an abstractNum's <w:tmpl> (or\n// <w:nsid>) value, so it cannot be targeted directly.\n//\n// The only API surface that touches numbering XML at all is the list\n// surface reachable from a paragraph (`paragraph.startNewList()`,\n// `paragraph.attachToList()`, `list.setLevel*()`, \u2026). Exercising any of\n// those rewrites the *entire* <w:numbering> part from scratch, which drops\n// <w:tmpl> (and <w:nsid>'s position, and w15:restartNumberingAfterBreak)\n// from *every* abstractNum in the document, not just the ten this diff\n// touches - i.e. it would overwrite far more than the requested change and\n// move the document further from the target, not closer. It also only\n// reaches templates already attached to body content, and none of\n// abstractNumId 1-10 are (no numbered/bulleted paragraph exists in this\n// document), so there is no live object to carry the new GUID even if one\n// were written.\n//\n// Given there is no Word.* API call that can express \"change this\n// abstractNum's tmpl GUID\" without that destructive side effect, this\n// script intentionally performs a no-op content-preserving read so the\n// document is left byte-for-byte as close to the (unreachable) target as\n// possible rather than corrupting unrelated numbering definitions.\n// ---------------------------------------------------------------------------\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\n// Confirm there is nothing in the document that actually consumes any of\n// the legacy numbering definitions (abstractNumId 1-10) the diff touches -\n// every paragraph uses ordinary paragraph styles, never list formatting,\n// so there is no ListFormat/List/ListLevel object reachable from body\n// content that could carry a <w:tmpl> value even indirectly.\nfor (const paragraph of body.paragraphs.items) {\n  paragraph.load(\"isListItem\");\n}\nawait context.sync();\n\n// No reachable Word.js property corresponds to w:tmpl - intentionally no\n// further action is taken so word/numbering.xml is left untouched.\n", "ps1": "# ---------------------------------------------------------------------------\n# Target change (see XML diff in the task description):\n#   word/numbering.xml - for the ten legacy <w:abstractNum> definitions\n#   (abstractNumId=\"1\"..\"10\", the unused \"FFFFFF7C\"..\"FFFFFF89\" nsid block\n#   inherited from the Pandoc/R-Markdown \"reference docx\" template) the\n#   cached <w:tmpl w:val=\"\u2026\"/> correlation GUID is swapped for a new,\n#   unrelated GUID. Nothing else in those definitions (nsid, lvl, numFmt,\n#   indentation, \u2026) changes, and no paragraph in the document body actually\n#   uses any of these numbering definitions (no numPr anywhere in Content).\n#\n# <w:tmpl> is Word's internal \"list template cache key\". It is not content,\n# not formatting, and is not surfaced by any property in the Word object\n# model - ListTemplates/ListTemplate/ListLevels/ListLevel expose Name,\n# NumberFormat, StartAt, Alignment, TabPosition, etc., but nothing maps to\n# an abstractNum's <w:tmpl> (or <w:nsid>) value, so it can't be targeted\n# directly. (Word's own ListTemplates collection here also only returns\n# generic placeholder levels for templates that aren't attached to any\n# Range/Paragraph - these ten never are - confirming there is no live\n# handle onto their stored XML to begin with.)\n#\n# The only things that touch numbering XML at all are operations that mint\n# or attach list formatting on a Paragraph/Range (e.g. ListFormat.ApplyListTemplate,\n# Range.ListFormat.ApplyNumberDefault, attaching a new ListTemplate, \u2026).\n# Exercising any of those rewrites the *entire* numbering part from scratch,\n# which drops <w:tmpl> (and relocates <w:nsid>, and drops\n# w15:restartNumberingAfterBreak) from *every* abstractNum in the document,\n# not just the ten this diff touches - i.e. it would overwrite far more than\n# the requested change and move the document further from the target, not\n# closer.\n#\n# Given there is no Word COM call that can express \"change this\n# abstractNum's tmpl GUID\" without that destructive side effect, this\n# script intentionally performs a no-op content-preserving read so the\n# document is left byte-for-byte as close to the (unreachable) target as\n# possible rather than corrupting unrelated numbering definitions.\n# ---------------------------------------------------------------------------\n\n$d = $word.ActiveDocument\n\n# Confirm there is nothing in the document that actually consumes any of\n# the legacy numbering definitions (abstractNumId 1-10) the diff touches -\n# walk the paragraphs/content and show none carry list formatting, so there\n# is no ListFormat/ListTemplate/ListLevel object reachable from body content\n# that could carry a <w:tmpl> value even indirectly.\n$range = $d.Range()\nforeach ($para in $d.Paragraphs) {\n    $null = $para.Range.ListFormat.ListType\n}\n$null = $range.Find.Execute(\"\")\n\n# No reachable Word COM property corresponds to w:tmpl - intentionally no\n# further action is taken so word/numbering.xml is left untouched.\n"}
